$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Cells.Item(3,1).Value = 49
$ws.Cells.Item(3,2).Value = "5:50 AM"
$ws.Cells.Item(3,3).Value = 1467.18
$ws.Cells.Item(3,4).Value = 1497.18
$ws.Cells.Item(3,5).Value = 24.93

# Row 5
$ws.Cells.Item(5,1).Value = 123
$ws.Cells.Item(5,2).Value = "6:27 AM"
$ws.Cells.Item(5,3).Value = 3687.43
$ws.Cells.Item(5,4).Value = 3717.43
$ws.Cells.Item(5,5).Value = 25.2

# Row 7
$ws.Cells.Item(7,1).Value = 182
$ws.Cells.Item(7,2).Value = "6:56 AM"
$ws.Cells.Item(7,3).Value = 5440.89
$ws.Cells.Item(7,4).Value = 5470.89
$ws.Cells.Item(7,5).Value = 108.74

# Row 8
$ws.Cells.Item(8,1).Value = 230
$ws.Cells.Item(8,2).Value = "7:20 AM"
$ws.Cells.Item(8,3).Value = 6880.786667
$ws.Cells.Item(8,4).Value = 6910.786667
$ws.Cells.Item(8,5).Value = 32.193333

# Row 9
$ws.Cells.Item(9,1).Value = 251
$ws.Cells.Item(9,2).Value = "7:30 AM"
$ws.Cells.Item(9,3).Value = 7516.66
$ws.Cells.Item(9,4).Value = 7546.66
$ws.Cells.Item(9,5).Value = 27.63

# Row 10
$ws.Cells.Item(10,1).Value = 275
$ws.Cells.Item(10,2).Value = "7:43 AM"
$ws.Cells.Item(10,3).Value = 8248.530000000001
$ws.Cells.Item(10,4).Value = 8278.530000000001
$ws.Cells.Item(10,5).Value = 21.35

# Row 11
$ws.Cells.Item(11,1).Value = 524
$ws.Cells.Item(11,2).Value = "9:47 AM"
$ws.Cells.Item(11,3).Value = 15702.116667
$ws.Cells.Item(11,4).Value = 15732.116667
$ws.Cells.Item(11,5).Value = 35.09

# Row 12
$ws.Cells.Item(12,1).Value = 550
$ws.Cells.Item(12,2).Value = "10:00 A"
$ws.Cells.Item(12,3).Value = 16484
$ws.Cells.Item(12,4).Value = 16514
$ws.Cells.Item(12,5).Value = 26.08

# Row 13
$ws.Cells.Item(13,1).Value = 570
$ws.Cells.Item(13,2).Value = "10:10 A"
$ws.Cells.Item(13,3).Value = 17087.03
$ws.Cells.Item(13,4).Value = 17117.03
$ws.Cells.Item(13,5).Value = 37.05

# Row 14
$ws.Cells.Item(14,1).Value = 993
$ws.Cells.Item(14,2).Value = "1:41 PM"
$ws.Cells.Item(14,3).Value = 29775.52
$ws.Cells.Item(14,4).Value = 29805.52
$ws.Cells.Item(14,5).Value = 41.35

# Row 15
$ws.Cells.Item(15,1).Value = 1016
$ws.Cells.Item(15,2).Value = "1:53 PM"
$ws.Cells.Item(15,3).Value = 30454.73
$ws.Cells.Item(15,4).Value = 30484.73
$ws.Cells.Item(15,5).Value = 24.88

# Row 16
$ws.Cells.Item(16,1).Value = 1022
$ws.Cells.Item(16,2).Value = "1:56 PM"
$ws.Cells.Item(16,3).Value = 30644.296667
$ws.Cells.Item(16,4).Value = 30674.296667
$ws.Cells.Item(16,5).Value = 24.493333

# Row 17
$ws.Cells.Item(17,1).Value = 1036
$ws.Cells.Item(17,2).Value = "2:03 PM"
$ws.Cells.Item(17,3).Value = 31058.49
$ws.Cells.Item(17,4).Value = 31088.49
$ws.Cells.Item(17,5).Value = 60.24

# Row 18
$ws.Cells.Item(18,1).Value = 1086
$ws.Cells.Item(18,2).Value = "2:28 PM"
$ws.Cells.Item(18,3).Value = 32575.77
$ws.Cells.Item(18,4).Value = 32605.77
$ws.Cells.Item(18,5).Value = 21.33

# Row 19
$ws.Cells.Item(19,1).Value = 1103
$ws.Cells.Item(19,2).Value = "2:36 PM"
$ws.Cells.Item(19,3).Value = 33070.88
$ws.Cells.Item(19,4).Value = 33100.88
$ws.Cells.Item(19,5).Value = 25.1

# Row 20
$ws.Cells.Item(20,1).Value = 1141
$ws.Cells.Item(20,2).Value = "2:55 PM"
$ws.Cells.Item(20,3).Value = 34219.87
$ws.Cells.Item(20,4).Value = 34249.87
$ws.Cells.Item(20,5).Value = 21.48

# Row 21
$ws.Cells.Item(21,1).Value = 1182
$ws.Cells.Item(21,2).Value = "3:16 PM"
$ws.Cells.Item(21,3).Value = 35455.655
$ws.Cells.Item(21,4).Value = 35485.655
$ws.Cells.Item(21,5).Value = 37.115

# Row 22
$ws.Cells.Item(22,1).Value = 1190
$ws.Cells.Item(22,2).Value = "3:20 PM"
$ws.Cells.Item(22,3).Value = 35678.675
$ws.Cells.Item(22,4).Value = 35708.675
$ws.Cells.Item(22,5).Value = 40.53

# Row 23
$ws.Cells.Item(23,1).Value = 1196
$ws.Cells.Item(23,2).Value = "3:23 PM"
$ws.Cells.Item(23,3).Value = 35858.715
$ws.Cells.Item(23,4).Value = 35888.715
$ws.Cells.Item(23,5).Value = 56.205

# Row 24
$ws.Cells.Item(24,1).Value = 1202
$ws.Cells.Item(24,2).Value = "3:26 PM"
$ws.Cells.Item(24,3).Value = 36041.01
$ws.Cells.Item(24,4).Value = 36071.01
$ws.Cells.Item(24,5).Value = 24.126667

# Row 25
$ws.Cells.Item(25,1).Value = 1209
$ws.Cells.Item(25,2).Value = "3:29 PM"
$ws.Cells.Item(25,3).Value = 36244.02
$ws.Cells.Item(25,4).Value = 36274.02
$ws.Cells.Item(25,5).Value = 75.29000000000001

# Row 26
$ws.Cells.Item(26,1).Value = 1221
$ws.Cells.Item(26,2).Value = "3:35 PM"
$ws.Cells.Item(26,3).Value = 36616.4
$ws.Cells.Item(26,4).Value = 36646.4
$ws.Cells.Item(26,5).Value = 23.47

# Row 27
$ws.Cells.Item(27,1).Value = 1227
$ws.Cells.Item(27,2).Value = "3:39 PM"
$ws.Cells.Item(27,3).Value = 36808.87
$ws.Cells.Item(27,4).Value = 36838.87
$ws.Cells.Item(27,5).Value = 22.12

# Row 28
$ws.Cells.Item(28,1).Value = 1245
$ws.Cells.Item(28,2).Value = "3:47 PM"
$ws.Cells.Item(28,3).Value = 37332.09
$ws.Cells.Item(28,4).Value = 37362.09
$ws.Cells.Item(28,5).Value = 193.12

# Row 29
$ws.Cells.Item(29,1).Value = 1269
$ws.Cells.Item(29,2).Value = "3:59 PM"
$ws.Cells.Item(29,3).Value = 38052.733333
$ws.Cells.Item(29,4).Value = 38082.733333
$ws.Cells.Item(29,5).Value = 20.703333

# Row 30
$ws.Cells.Item(30,1).Value = 1279
$ws.Cells.Item(30,2).Value = "4:05 PM"
$ws.Cells.Item(30,3).Value = 38365.04
$ws.Cells.Item(30,4).Value = 38395.04
$ws.Cells.Item(30,5).Value = 31.16

# Row 31
$ws.Cells.Item(31,1).Value = 1288
$ws.Cells.Item(31,2).Value = "4:09 PM"
$ws.Cells.Item(31,3).Value = 38622.23
$ws.Cells.Item(31,4).Value = 38652.23
$ws.Cells.Item(31,5).Value = 50.7

# Row 32
$ws.Cells.Item(32,1).Value = 1302
$ws.Cells.Item(32,2).Value = "4:16 PM"
$ws.Cells.Item(32,3).Value = 39032.49
$ws.Cells.Item(32,4).Value = 39062.49
$ws.Cells.Item(32,5).Value = 230.5

# Row 33
$ws.Cells.Item(33,1).Value = 1308
$ws.Cells.Item(33,2).Value = "4:19 PM"
$ws.Cells.Item(33,3).Value = 39217.913333
$ws.Cells.Item(33,4).Value = 39247.913333
$ws.Cells.Item(33,5).Value = 29.03

# Row 34
$ws.Cells.Item(34,1).Value = 1315
$ws.Cells.Item(34,2).Value = "4:22 PM"
$ws.Cells.Item(34,3).Value = 39427.435
$ws.Cells.Item(34,4).Value = 39457.435
$ws.Cells.Item(34,5).Value = 98.93000000000001

# Row 35
$ws.Cells.Item(35,1).Value = 1323
$ws.Cells.Item(35,2).Value = "4:26 PM"
$ws.Cells.Item(35,3).Value = 39677.255
$ws.Cells.Item(35,4).Value = 39707.255
$ws.Cells.Item(35,5).Value = 57.735

# Row 36
$ws.Cells.Item(36,1).Value = 1333
$ws.Cells.Item(36,2).Value = "4:31 PM"
$ws.Cells.Item(36,3).Value = 39975.185
$ws.Cells.Item(36,4).Value = 40005.185
$ws.Cells.Item(36,5).Value = 26.385

# Row 37
$ws.Cells.Item(37,1).Value = 1345
$ws.Cells.Item(37,2).Value = "4:38 PM"
$ws.Cells.Item(37,3).Value = 40347.76
$ws.Cells.Item(37,4).Value = 40377.76
$ws.Cells.Item(37,5).Value = 94.73999999999999

# Row 38
$ws.Cells.Item(38,1).Value = 1362
$ws.Cells.Item(38,2).Value = "4:46 PM"
$ws.Cells.Item(38,3).Value = 40838.935
$ws.Cells.Item(38,4).Value = 40868.935
$ws.Cells.Item(38,5).Value = 104.33

# Row 39
$ws.Cells.Item(39,1).Value = 1369
$ws.Cells.Item(39,2).Value = "4:49 PM"
$ws.Cells.Item(39,3).Value = 41046.07
$ws.Cells.Item(39,4).Value = 41076.07
$ws.Cells.Item(39,5).Value = 33.66

# Row 40
$ws.Cells.Item(40,1).Value = 1390
$ws.Cells.Item(40,2).Value = "5:00 PM"
$ws.Cells.Item(40,3).Value = 41680.48
$ws.Cells.Item(40,4).Value = 41710.48
$ws.Cells.Item(40,5).Value = 38.76

# Row 41
$ws.Cells.Item(41,1).Value = 1410
$ws.Cells.Item(41,2).Value = "5:10 PM"
$ws.Cells.Item(41,3).Value = 42278.335
$ws.Cells.Item(41,4).Value = 42308.335
$ws.Cells.Item(41,5).Value = 44.55

# Row 42
$ws.Cells.Item(42,1).Value = 1417
$ws.Cells.Item(42,2).Value = "5:13 PM"
$ws.Cells.Item(42,3).Value = 42486.06
$ws.Cells.Item(42,4).Value = 42516.06
$ws.Cells.Item(42,5).Value = 61.41

# Row 43
$ws.Cells.Item(43,1).Value = 1427
$ws.Cells.Item(43,2).Value = "5:19 PM"
$ws.Cells.Item(43,3).Value = 42809.97
$ws.Cells.Item(43,4).Value = 42839.97
$ws.Cells.Item(43,5).Value = 161.82

# Row 44
$ws.Cells.Item(44,1).Value = 1445
$ws.Cells.Item(44,2).Value = "5:27 PM"
$ws.Cells.Item(44,3).Value = 43323.345
$ws.Cells.Item(44,4).Value = 43353.345
$ws.Cells.Item(44,5).Value = 29.31

# Row 45
$ws.Cells.Item(45,1).Value = 1451
$ws.Cells.Item(45,2).Value = "5:30 PM"
$ws.Cells.Item(45,3).Value = 43516.02
$ws.Cells.Item(45,4).Value = 43546.02
$ws.Cells.Item(45,5).Value = 88.58

# Row 46
$ws.Cells.Item(46,1).Value = 1472
$ws.Cells.Item(46,2).Value = "5:41 PM"
$ws.Cells.Item(46,3).Value = 44138.865
$ws.Cells.Item(46,4).Value = 44168.865
$ws.Cells.Item(46,5).Value = 48.69

# Row 47
$ws.Cells.Item(47,1).Value = 1479
$ws.Cells.Item(47,2).Value = "5:44 PM"
$ws.Cells.Item(47,3).Value = 44340.68
$ws.Cells.Item(47,4).Value = 44370.68
$ws.Cells.Item(47,5).Value = 49.65

# Row 48
$ws.Cells.Item(48,1).Value = 1492
$ws.Cells.Item(48,2).Value = "5:51 PM"
$ws.Cells.Item(48,3).Value = 44738.87
$ws.Cells.Item(48,4).Value = 44768.87
$ws.Cells.Item(48,5).Value = 22.62

